# Edit script: restructure condition table, adding stimRamp/rampLin columns,
# renaming dnDivUp -> stepUpMulti, adjusting condition values, and updating the view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before the old "postStimBlankT" column (J),
#    shifting postStimBlankT..nRevs (J:S) to L:U.
$ws.Columns("J:K").Insert()

# 2. Rename the (now shifted) dnDivUp header to stepUpMulti, and bump its
#    value from 1 to 4 for every condition row.
$ws.Range("T1").Value = "stepUpMulti"
$ws.Range("T2:T12").Value = 4

# 3. New column headers for the inserted columns.
$ws.Range("J1").Value = "stimRamp"
$ws.Range("K1").Value = "rampLin"

# 4. New column values (stimRamp = 1, rampLin = 1) for all data rows.
$ws.Range("J2:J12").Value = 1
$ws.Range("K2:K12").Value = 1

# 5. singlCont (column D) drops to 0 for the first contrast pair (rows 2-3).
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0

# 6. jitTmax (column E) changes from 500 to 250 for every condition row.
$ws.Range("E2:E12").Value = 250

# 7. postStimBlankT (now column L) changes from 500 to 250 for every row.
$ws.Range("L2:L12").Value = 250

# 8. odtLoc (now column N) changes from 6 to 1 for the "opposite" rows.
$ws.Range("N3").Value = 1
$ws.Range("N5").Value = 1
$ws.Range("N7").Value = 1
$ws.Range("N9").Value = 1
$ws.Range("N11").Value = 1

# 9. stairUp (now column R) changes from 1 to 2 for every row.
$ws.Range("R2:R12").Value = 2

# 10. Update the view: zoom to 130% and select K2.
$aw = $excel.ActiveWindow
$aw.Zoom = 130
$ws.Range("K2").Select()
